$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44355
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 18000
$ws.Range("S2").Value = 1000

# Row 3
$ws.Range("D3").Value = 44348
$ws.Range("M3").Value = 200

# Row 4
$ws.Range("D4").Value = 44326
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 1111

# Row 7
$ws.Range("D7").Value = 44340
$ws.Range("M7").Value = 230
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("Q7").Value = "$/caja 18 kilos granel"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 1111
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 44291
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 150
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("Q8").Value = "$/caja 15 kilos granel"
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 800
$ws.Range("T8").Value = 15

# Row 9
$ws.Range("D9").Value = 44342
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 1111

# Row 10
$ws.Range("D10").Value = 44319
$ws.Range("M10").Value = 120

# Row 13
$ws.Range("D13").Value = 44328
$ws.Range("L13").Value = "Especial"
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range("Q13").Value = "$/caja 18 kilos granel"
$ws.Range("R13").Value = "Provincia de Limarí"
$ws.Range("S13").Value = 1111
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = 44294
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("Q14").Value = "$/caja 15 kilos granel"
$ws.Range("R14").Value = "Región Metropolitana"
$ws.Range("S14").Value = 800
$ws.Range("T14").Value = 15

# Row 15
$ws.Range("D15").Value = 44354
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 18000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 18000
$ws.Range("S15").Value = 1000
